$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent values for case with 380 kV (rows 2-25, i.e. data rows 0-23)
# Column B
$ws.Range("B2").Value = 19.62974171293408
$ws.Range("B3").Value = 19.39771230713015
$ws.Range("B4").Value = 19.26257364858353
$ws.Range("B5").Value = 19.20941934262618
$ws.Range("B6").Value = 19.20071085665069
$ws.Range("B7").Value = 19.26184894114633
$ws.Range("B8").Value = 19.54826225844449
$ws.Range("B9").Value = 20.16463181465948
$ws.Range("B10").Value = 20.64598627138448
$ws.Range("B11").Value = 20.87008020179707
$ws.Range("B12").Value = 20.95558421132854
$ws.Range("B13").Value = 20.93714219377285
$ws.Range("B14").Value = 20.87710227862311
$ws.Range("B15").Value = 20.84040729711566
$ws.Range("B16").Value = 20.63143641370438
$ws.Range("B17").Value = 20.5044861240043
$ws.Range("B18").Value = 20.43195545097134
$ws.Range("B19").Value = 20.40748439909499
$ws.Range("B20").Value = 20.51795038483303
$ws.Range("B21").Value = 20.8947206991028
$ws.Range("B22").Value = 21.14467611788595
$ws.Range("B23").Value = 21.01096071218062
$ws.Range("B24").Value = 20.5118617691013
$ws.Range("B25").Value = 19.99252008330669

# Column C
$ws.Range("C2").Value = 11.90515418833905
$ws.Range("C3").Value = 11.71739948170048
$ws.Range("C4").Value = 11.60724474835997
$ws.Range("C5").Value = 11.56371910496087
$ws.Range("C6").Value = 11.55657613311774
$ws.Range("C7").Value = 11.60665213131123
$ws.Range("C8").Value = 11.83939763723914
$ws.Range("C9").Value = 12.33308491229793
$ws.Range("C10").Value = 12.71381497132857
$ws.Range("C11").Value = 12.88992159134418
$ws.Range("C12").Value = 12.95694422062275
$ws.Range("C13").Value = 12.94249611465856
$ws.Range("C14").Value = 12.89542934394477
$ws.Range("C15").Value = 12.86664074978992
$ws.Range("C16").Value = 12.70235736237088
$ws.Range("C17").Value = 12.60226018931339
$ws.Range("C18").Value = 12.54496581472114
$ws.Range("C19").Value = 12.52561740052046
$ws.Range("C20").Value = 12.61288745439396
$ws.Range("C21").Value = 12.90924556362616
$ws.Range("C22").Value = 13.10484741288149
$ws.Range("C23").Value = 13.0003029076297
$ws.Range("C24").Value = 12.60808207399127
$ws.Range("C25").Value = 12.19602102589494

# Column D
$ws.Range("D2").Value = 5.273071938883459
$ws.Range("D3").Value = 5.216086276269216
$ws.Range("D4").Value = 5.181422748973388
$ws.Range("D5").Value = 5.167386590410773
$ws.Range("D6").Value = 5.165061517764067
$ws.Range("D7").Value = 5.181233080257248
$ws.Range("D8").Value = 5.253360266723604
$ws.Range("D9").Value = 5.397043490281657
$ws.Range("D10").Value = 5.503501399747526
$ws.Range("D11").Value = 5.552022358177561
$ws.Range("D12").Value = 5.570400427781181
$ws.Range("D13").Value = 5.566442338025453
$ws.Range("D14").Value = 5.553534300203744
$ws.Range("D15").Value = 5.545628036081167
$ws.Range("D16").Value = 5.500331529110892
$ws.Range("D17").Value = 5.472560866529843
$ws.Range("D18").Value = 5.456597189552176
$ws.Range("D19").Value = 5.451194032743079
$ws.Range("D20").Value = 5.475516207091133
$ws.Range("D21").Value = 5.55732566130694
$ws.Range("D22").Value = 5.610813356395066
$ws.Range("D23").Value = 5.582267092636058
$ws.Range("D24").Value = 5.474180090557915
$ws.Range("D25").Value = 5.357981458067466

# Column E
$ws.Range("E2").Value = 10.73751006893296
$ws.Range("E3").Value = 10.77371095808514
$ws.Range("E4").Value = 10.79727827854848
$ws.Range("E5").Value = 10.80722020616768
$ws.Range("E6").Value = 10.80889150686934
$ws.Range("E7").Value = 10.79741098856526
$ws.Range("E8").Value = 10.74971474115319
$ws.Range("E9").Value = 10.66676207137593
$ws.Range("E10").Value = 10.61219537240914
$ws.Range("E11").Value = 10.58874159962467
$ws.Range("E12").Value = 10.58005594899964
$ws.Range("E13").Value = 10.58191786851938
$ws.Range("E14").Value = 10.58802310816462
$ws.Range("E15").Value = 10.59178820925135
$ws.Range("E16").Value = 10.61375558635303
$ws.Range("E17").Value = 10.62758170934053
$ws.Range("E18").Value = 10.63566303916502
$ws.Range("E19").Value = 10.63842141136882
$ws.Range("E20").Value = 10.62609656101013
$ws.Range("E21").Value = 10.58622454617681
$ws.Range("E22").Value = 10.56130657257316
$ws.Range("E23").Value = 10.5745017389685
$ws.Range("E24").Value = 10.62676758375843
$ws.Range("E25").Value = 10.68807771864422

# Column F
$ws.Range("F2").Value = 58.40536167705665
$ws.Range("F3").Value = 57.47542236806509
$ws.Range("F4").Value = 56.90513245647881
$ws.Range("F5").Value = 56.67312497238495
$ws.Range("F6").Value = 56.63463015741009
$ws.Range("F7").Value = 56.90200165818151
$ws.Range("F8").Value = 58.08470302662953
$ws.Range("F9").Value = 60.40010644014433
$ws.Range("F10").Value = 62.08648415599741
$ws.Range("F11").Value = 62.84791786180784
$ws.Range("F12").Value = 63.13522761613337
$ws.Range("F13").Value = 63.07339929687621
$ws.Range("F14").Value = 62.87157678723235
$ws.Range("F15").Value = 62.74781467270179
$ws.Range("F16").Value = 62.03658987126168
$ws.Range("F17").Value = 61.59866427041018
$ws.Range("F18").Value = 61.34625577150091
$ws.Range("F19").Value = 61.26071083113141
$ws.Range("F20").Value = 61.64533809227858
$ws.Range("F21").Value = 62.9308864863202
$ws.Range("F22").Value = 63.76498006604591
$ws.Range("F23").Value = 63.32043217726871
$ws.Range("F24").Value = 61.62423883264157
$ws.Range("F25").Value = 59.77543653153543

# Column H
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("H25").Value = 7.344005520526261

# Column J
$ws.Range("J2").Value = 10.20398338144533
$ws.Range("J3").Value = 10.20837293539282
$ws.Range("J4").Value = 10.21193096383938
$ws.Range("J5").Value = 10.21359756955083
$ws.Range("J6").Value = 10.21388738643812
$ws.Range("J7").Value = 10.21195256331966
$ws.Range("J8").Value = 10.20531765310049
$ws.Range("J9").Value = 10.19916581532182
$ws.Range("J10").Value = 10.19884609355878
$ws.Range("J11").Value = 10.19961630304799
$ws.Range("J12").Value = 10.20003986407436
$ws.Range("J13").Value = 10.19994277294555
$ws.Range("J14").Value = 10.19964850496848
$ws.Range("J15").Value = 10.19948544098091
$ws.Range("J16").Value = 10.19881420395936
$ws.Range("J17").Value = 10.19863711269807
$ws.Range("J18").Value = 10.19862142771259
$ws.Range("J19").Value = 10.19863091002314
$ws.Range("J20").Value = 10.19864704403269
$ws.Range("J21").Value = 10.19973135717269
$ws.Range("J22").Value = 10.20120892047465
$ws.Range("J23").Value = 10.20034989225893
$ws.Range("J24").Value = 10.19864228582691
$ws.Range("J25").Value = 10.2000937335388

# Column M
$ws.Range("M2").Value = 19.38323440442745
$ws.Range("M3").Value = 19.39223646570851
$ws.Range("M4").Value = 19.40398803934865
$ws.Range("M5").Value = 19.41033945220009
$ws.Range("M6").Value = 19.41148837811554
$ws.Range("M7").Value = 19.40406737413618
$ws.Range("M8").Value = 19.3850452064388
$ws.Range("M9").Value = 19.39720184337342
$ws.Range("M10").Value = 19.43631965603098
$ws.Range("M11").Value = 19.46065497148231
$ws.Range("M12").Value = 19.47080745234929
$ws.Range("M13").Value = 19.4685793147398
$ws.Range("M14").Value = 19.46147146159652
$ws.Range("M15").Value = 19.45723963412994
$ws.Range("M16").Value = 19.43486064130146
$ws.Range("M17").Value = 19.42280503746828
$ws.Range("M18").Value = 19.41648685714536
$ws.Range("M19").Value = 19.41445348995012
$ws.Range("M20").Value = 19.42402465561604
$ws.Range("M21").Value = 19.46353380604304
$ws.Range("M22").Value = 19.49481587292924
$ws.Range("M23").Value = 19.47762177717754
$ws.Range("M24").Value = 19.42347135744643
$ws.Range("M25").Value = 19.38861031543704
